# The six "bar" rectangles on slide 1 each grow taller by 64008 EMU
# (5.04 pt), with their Top/Left/Width left untouched (so each bar's
# bottom edge moves further down while its top edge stays put).
#
# Shape.Height is exposed in points and, in this host, is round-tripped
# through a single-precision float before being converted back to EMU
# (with truncation), so naive "Height = Height + deltaPoints" arithmetic
# can land 1 EMU short. To guarantee the exact target EMU value we
# compute the desired EMU explicitly and then nudge the assigned point
# value up/down in tiny steps until the shape reports back exactly that
# EMU amount.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerPoint = 12700.0
$deltaEmu = 64008   # +5.04 pt, applied to every bar's height

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)

    $curEmu = [math]::Round($sh.Height * $emuPerPoint)
    $targetEmu = $curEmu + $deltaEmu
    $targetPt = $targetEmu / $emuPerPoint

    $sh.Height = $targetPt
    $readEmu = [math]::Round($sh.Height * $emuPerPoint)

    $tries = 0
    while ($readEmu -ne $targetEmu -and $tries -lt 50) {
        if ($readEmu -lt $targetEmu) {
            $targetPt = $targetPt + 0.00001
        } else {
            $targetPt = $targetPt - 0.00001
        }
        $sh.Height = $targetPt
        $readEmu = [math]::Round($sh.Height * $emuPerPoint)
        $tries++
    }
}
